# Update countries & provincias Spain
# Refresh the COVID dashboard ("Pais" sheet) with the newer data snapshot
# (15 de Mayo de 2020 a las 01:05) and re-rank the few countries whose
# "Casos totales" (col B) overtook their neighbour in the table, which is
# sorted descending by column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / timestamp banner (A1)
$ws.Range("A1").Value = 'Datos actualizados a 15 de Mayo de 2020 a las 01:05'

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1454448
$ws.Range("C4").Value = 24100
$ws.Range("D4").Value = 316305
$ws.Range("E4").Value = 1051339
$ws.Range("G4").Value = 1607
$ws.Range("H4").Value = 86804

# Brasil (row 9)
$ws.Range("B9").Value = 202918
$ws.Range("C9").Value = 13761
$ws.Range("D9").Value = 79479
$ws.Range("E9").Value = 109446
$ws.Range("G9").Value = 835
$ws.Range("H9").Value = 13993

# Canada (row 17)
$ws.Range("B17").Value = 73400
$ws.Range("C17").Value = 1122
$ws.Range("D17").Value = 36091
$ws.Range("E17").Value = 31837
$ws.Range("G17").Value = 170
$ws.Range("H17").Value = 5472

# Chequia (row 51)
$ws.Range("B51").Value = 8351
$ws.Range("C51").Value = 82
$ws.Range("D51").Value = 5241
$ws.Range("E51").Value = 2817
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 293

# Noruega (row 52)
$ws.Range("B52").Value = 8196
$ws.Range("C52").Value = 21
$ws.Range("E52").Value = 7932

# Argentina overtakes Australia -> swap rows 53/54
$ws.Range("A53").Value = 'Argentina'
$ws.Range("B53").Value = 7134
$ws.Range("C53").Value = 255
$ws.Range("D53").Value = 2385
$ws.Range("E53").Value = 4396
$ws.Range("F53").Value = 170
$ws.Range("G53").Value = 24
$ws.Range("H53").Value = 353

$ws.Range("A54").Value = 'Australia'
$ws.Range("B54").Value = 6989
$ws.Range("C54").Value = 9
$ws.Range("D54").Value = 6301
$ws.Range("E54").Value = 590
$ws.Range("F54").Value = 18
$ws.Range("G54").Value = 0
$ws.Range("H54").Value = 98

# Barein overtakes Finlandia -> swap rows 58/59
$ws.Range("A58").Value = 'Barein'
$ws.Range("B58").Value = 6198
$ws.Range("C58").Value = 382
$ws.Range("D58").Value = 2353
$ws.Range("E58").Value = 3835
$ws.Range("F58").Value = 4
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 10

$ws.Range("A59").Value = 'Finlandia'
$ws.Range("B59").Value = 6145
$ws.Range("C59").Value = 91
$ws.Range("D59").Value = 4300
$ws.Range("E59").Value = 1558
$ws.Range("F59").Value = 33
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = 287

# Maldivas (row 102)
$ws.Range("B102").Value = 982
$ws.Range("C102").Value = 27
$ws.Range("D102").Value = 45
$ws.Range("E102").Value = 933

# Uruguay (row 117)
$ws.Range("B117").Value = 724
$ws.Range("C117").Value = 5
$ws.Range("D117").Value = 547
$ws.Range("E117").Value = 158
$ws.Range("F117").Value = 7

# Guayana Francesa overtakes Uganda and Guadalupe -> rows 153/154/155 shift
$ws.Range("A153").Value = 'Guayana Francesa'
$ws.Range("B153").Value = 164
$ws.Range("C153").Value = 11
$ws.Range("D153").Value = 124
$ws.Range("E153").Value = 39
$ws.Range("H153").Value = 1

$ws.Range("A154").Value = 'Uganda'
$ws.Range("B154").Value = 160
$ws.Range("C154").Value = 21
$ws.Range("D154").Value = 63
$ws.Range("E154").Value = 97
$ws.Range("F154").Value = 0
$ws.Range("H154").Value = 0

$ws.Range("A155").Value = 'Guadalupe'
$ws.Range("B155").Value = 155
$ws.Range("D155").Value = 109
$ws.Range("E155").Value = 33
$ws.Range("F155").Value = 4
$ws.Range("H155").Value = 13

# Guyana (row 163)
$ws.Range("D163").Value = 42
$ws.Range("E163").Value = 61

# Bahamas overtakes Monaco -> swap rows 166/167
$ws.Range("A166").Value = 'Bahamas'
$ws.Range("C166").Value = 2
$ws.Range("D166").Value = 41
$ws.Range("E166").Value = 44
$ws.Range("H166").Value = 11

$ws.Range("A167").Value = 'Monaco'
$ws.Range("B167").Value = 96
$ws.Range("D167").Value = 87
$ws.Range("E167").Value = 5
$ws.Range("H167").Value = 4

# Polinesia Francesa (row 175)
$ws.Range("D175").Value = 59
$ws.Range("E175").Value = 1
$ws.Range("F175").Value = 0

# Mauritania (row 190)
$ws.Range("D190").Value = 7
$ws.Range("E190").Value = 11

# Comoras (row 207)
$ws.Range("D207").Value = 3
$ws.Range("E207").Value = 7
